$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as text so values like "1.00" or "0.0710" are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.926.03"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "3.530.93"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "588.75"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").Value = "177.40"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.528.43"
$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D11").Value = "6.95"
$ws.Range("E11").Value = "  -1.66%  "

$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "4.135.39"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").Value = "30.66"
$ws.Range("E14").Value = "  -3.93%  "

$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("D16").Value = "66.897.66"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").Value = "3.524.69"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").Value = "383.17"
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("D22").Value = "7.89"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "5.75"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("D26").Value = "71.78"
$ws.Range("E26").Value = "  -2.79%  "

$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -3.59%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").Value = "24.62"
$ws.Range("E31").Value = "  +4.75%  "

$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  -2.27%  "

$ws.Range("D33").Value = "2.03"
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("E35").Value = "  -1.14%  "

$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "1.58"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "29.69"
$ws.Range("E38").Value = "  +13.51%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "159.31"
$ws.Range("E39").Value = "  -3.22%  "

$ws.Range("D40").Value = "0.895"
$ws.Range("E40").Value = "  +3.03%  "

$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "6.64"
$ws.Range("E42").Value = "  -2.50%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.54"
$ws.Range("E43").Value = "  -2.34%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").Value = "  -5.76%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0710"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.729.41"
$ws.Range("E46").Value = "  -3.64%  "

$ws.Range("D47").Value = "25.55"
$ws.Range("E47").Value = "  -5.50%  "

$ws.Range("D48").Value = "40.73"
$ws.Range("E48").Value = "  -2.12%  "

$ws.Range("D49").Value = "0.0300"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("D50").Value = "327.79"
$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("E51").Value = "  -2.10%  "
